$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 309; this shifts existing rows 309-332
# down to 310-333, preserving all of their data/values untouched.
$ws.Rows.Item(309).Insert()

# Populate the newly inserted row 309 with the new data record.
$ws.Cells.Item(309, 1).Value = 10
$ws.Cells.Item(309, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(309, 3).Value = "La Araucanía"
$ws.Cells.Item(309, 4).NumberFormat = $ws.Cells.Item(310, 4).NumberFormat
$ws.Cells.Item(309, 4).Value = 44578
$ws.Cells.Item(309, 5).Value = 9
$ws.Cells.Item(309, 6).Value = 100114014
$ws.Cells.Item(309, 7).Value = "Betarraga"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 50
$ws.Cells.Item(309, 11).Value = 8000
$ws.Cells.Item(309, 12).Value = 8000
$ws.Cells.Item(309, 13).Value = 8000
$ws.Cells.Item(309, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(309, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(309, 16).Value = 667
$ws.Cells.Item(309, 17).Value = 12
$ws.Cells.Item(309, 18).Value = "Hortaliza"
